$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H28").Value = 2177.5715
$ws.Range("I28").Value = 2373.8333
$ws.Range("K28").Value = 2373.8333
$ws.Range("M28").Value = -1888.8333
$ws.Range("H70").Value = 2224
$ws.Range("I70").Value = 1999
$ws.Range("J70").Value = 2899
$ws.Range("K70").Value = 5997
$ws.Range("L70").Value = 8697
$ws.Range("M70").Value = -5727
$ws.Range("N70").Value = -9237
$ws.Range("H73").Value = 2224
$ws.Range("I73").Value = 1999
$ws.Range("J73").Value = 2899
$ws.Range("K73").Value = 5997
$ws.Range("L73").Value = 8697
$ws.Range("M73").Value = -5061
$ws.Range("N73").Value = -10569
$ws.Range("H76").Value = 900
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 900
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 900
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -1530
$ws.Range("H79").Value = 900
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 900
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 900
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -3084
$ws.Range("H80").Value = 1888.3334
$ws.Range("I80").Value = 948.75
$ws.Range("J80").Value = 2640
$ws.Range("K80").Value = 2846.25
$ws.Range("L80").Value = 7920
$ws.Range("M80").Value = -1848.25
$ws.Range("N80").Value = -9916
$ws.Range("H83").Value = 1888.3334
$ws.Range("I83").Value = 948.75
$ws.Range("J83").Value = 2640
$ws.Range("K83").Value = 8538.75
$ws.Range("L83").Value = 23760
$ws.Range("M83").Value = -3546.75
$ws.Range("N83").Value = -33744
$ws.Range("H138").Value = 4951.294
$ws.Range("I138").Value = 4108.1113
$ws.Range("K138").Value = 12324.3339
$ws.Range("M138").Value = -7184.333899999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1049.125
$ws.Range("I97").Value = 948.8333
$ws.Range("K97").Value = 948.8333
$ws.Range("M97").Value = -452.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1821.3334
$ws.Range("I20").Value = 200
$ws.Range("J20").Value = 2632
$ws.Range("K20").Value = 200
$ws.Range("L20").Value = 2632
$ws.Range("M20").Value = 47
$ws.Range("N20").Value = -3126
$ws.Range("H134").Value = 791.7
$ws.Range("I134").Value = 791.7
$ws.Range("K134").Value = 2375.1
$ws.Range("M134").Value = 159.8999999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2436.6128
$ws.Range("I31").Value = 1303.1428
$ws.Range("J31").Value = 2767.2083
$ws.Range("K31").Value = 1303.1428
$ws.Range("L31").Value = 2767.2083
$ws.Range("M31").Value = -1008.1428
$ws.Range("N31").Value = -3357.2083
$ws.Range("H34").Value = 2436.6128
$ws.Range("I34").Value = 1303.1428
$ws.Range("J34").Value = 2767.2083
$ws.Range("K34").Value = 1303.1428
$ws.Range("L34").Value = 2767.2083
$ws.Range("M34").Value = -1101.1428
$ws.Range("N34").Value = -3171.2083
$ws.Range("H58").Value = 3041.625
$ws.Range("I58").Value = 1552.5
$ws.Range("K58").Value = 1552.5
$ws.Range("M58").Value = -1349.5
$ws.Range("H136").Value = 3041.625
$ws.Range("I136").Value = 1552.5
$ws.Range("K136").Value = 4657.5
$ws.Range("M136").Value = -2107.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 201080.27
$ws.Range("I2").Value = 366721.34
$ws.Range("J2").Value = 138964.88
$ws.Range("K2").Value = 2200328.04
$ws.Range("L2").Value = 833789.28
$ws.Range("M2").Value = -2200215.04
$ws.Range("N2").Value = -834015.28
$ws.Range("H11").Value = 20853712
$ws.Range("I11").Value = 20853712
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 62561136
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -62560996
$ws.Range("N11").ClearContents()
$ws.Range("H22").Value = 2236.2666
$ws.Range("J22").Value = 2253.1428
$ws.Range("L22").Value = 6759.428400000001
$ws.Range("N22").Value = -7097.428400000001
$ws.Range("H23").Value = 347.125
$ws.Range("I23").Value = 98.5
$ws.Range("J23").Value = 430
$ws.Range("K23").Value = 295.5
$ws.Range("L23").Value = 1290
$ws.Range("M23").Value = -60.5
$ws.Range("N23").Value = -1760
$ws.Range("H26").Value = 1374.75
$ws.Range("I26").Value = 1299.6666
$ws.Range("J26").Value = 1600
$ws.Range("K26").Value = 3898.9998
$ws.Range("L26").Value = 4800
$ws.Range("M26").Value = -3610.9998
$ws.Range("N26").Value = -5376
$ws.Range("H27").Value = 2236.2666
$ws.Range("J27").Value = 2253.1428
$ws.Range("L27").Value = 6759.428400000001
$ws.Range("N27").Value = -6963.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2466
$ws.Range("I70").Value = 2224.5
$ws.Range("J70").Value = 2949
$ws.Range("K70").Value = 2224.5
$ws.Range("L70").Value = 2949
$ws.Range("M70").Value = -1954.5
$ws.Range("N70").Value = -3489
$ws.Range("H73").Value = 2466
$ws.Range("I73").Value = 2224.5
$ws.Range("J73").Value = 2949
$ws.Range("K73").Value = 2224.5
$ws.Range("L73").Value = 2949
$ws.Range("M73").Value = -1288.5
$ws.Range("N73").Value = -4821
$ws.Range("H107").Value = 1958.5
$ws.Range("I107").Value = 91.5
$ws.Range("J107").Value = 3825.5
$ws.Range("K107").Value = 91.5
$ws.Range("L107").Value = 3825.5
$ws.Range("M107").Value = 1828.5
$ws.Range("N107").Value = -7665.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3568.6
$ws.Range("I46").Value = 1937.2
$ws.Range("J46").Value = 5200
$ws.Range("K46").Value = 1937.2
$ws.Range("L46").Value = 5200
$ws.Range("M46").Value = -1749.2
$ws.Range("N46").Value = -5576
$ws.Range("H61").Value = 2658.3333
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2658.3333
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 2658.3333
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -3062.3333
$ws.Range("H100").Value = 1268.0714
$ws.Range("I100").Value = 1288.8462
$ws.Range("J100").Value = 998
$ws.Range("K100").Value = 1288.8462
$ws.Range("L100").Value = 998
$ws.Range("M100").Value = -747.8462
$ws.Range("N100").Value = -2080
$ws.Range("H113").Value = 2658.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2658.3333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2658.3333
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6998.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 33335232
$ws.Range("I100").Value = 50001450
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 100002900
$ws.Range("L100").Value = 5600
$ws.Range("M100").Value = -100002359
$ws.Range("N100").Value = -6682
$ws.Range("H113").Value = 669.4
$ws.Range("I113").Value = 658
$ws.Range("J113").Value = 686.5
$ws.Range("K113").Value = 1974
$ws.Range("L113").Value = 2059.5
$ws.Range("M113").Value = 196
$ws.Range("N113").Value = -6399.5
